# Replace the descriptive vehicle-class labels in the "data" sheet with
# their short code equivalents, matching the shared-string text swap:
#   "Passenger cars"            -> "NOPC"
#   "Combined vehicles"         -> "NOCV"
#   "Tractors etc"              -> "NOTR"
#   "Special purpose vehicles"  -> "NOSPV"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$map = @{
    "Passenger cars"           = "NOPC"
    "Combined vehicles"        = "NOCV"
    "Tractors etc"             = "NOTR"
    "Special purpose vehicles" = "NOSPV"
}

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 6)  # column F = vehicle_class
    $val = $cell.Value()
    if ($map.ContainsKey($val)) {
        $cell.Value = $map[$val]
    }
}
